$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the per-attribute rows (name/type/abilities/power-toughness) for
# each token into a single Python-tuple-style string in column A, rows 2-4,
# then drop the now-unused rows 5-14.
$ws.Range("A2").Value = "('Bat', ['Token Creature — Bat', 'Flying', '1/1'])"
$ws.Range("A3").Value = "('Beast', ['Token Creature — Beast', '4/4'])"
$ws.Range("A4").Value = "('Elephant', ['Token Creature — Elephant', '3/3'])"

$ws.Range("A5:A14").ClearContents()
